$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474758173/gFRHQATQ91-camera-video-segment-driver-1749474760673.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSI2GHGTCK%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T220004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIE660DgdFCi8i8bvMivDROKxD%2ByCQRzirgbx79lzYgn3AiAyRf7fqefgxSZENBHtm9Aw%2FgI06qdyTqqpz4idAdpiLCrmAwjF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0qByG9j973ShKj46KroD61BMkswDVOuyWr%2F39rq%2F15a2nLb9X25cXbmqw525F0vEhltm3owMssOpJ8XG1g8fzPysBM%2FV8DLxsUIGLquTJnrJTQo5Tp3OpsWS2vuWz0LGciUFUlCibXMA10D0KSqlpjwHNY7hw5WwRRrZrkgY4NLes1xcPXZIl9BtLavqr0Xnb5VpSFVYIkfLaOeISCUymmd1F4PmsEPADaLjm9ZGSv9s5%2BupvKbNMkzTAevyWRTlcDwXQ7O8AkQNphRQy5kKvP2wGkM29t%2FzAFqXq85D%2FDgx%2F1bMS4PMbALbXh85jfqYC7pqCeu63rpwsHdsrywwsegg0QnXznhlAgFKzVtORnt6xSG6TPuvy3DP1wphmT7sXLxjx2kSk1wM9WdQBwr%2BbmZr10hR0dj0Tmv5kEXOVT4exoZK8VXihIepJBcSu7va%2BECIrOlyo1WtQE52mocX0TG45W%2B9MNB5J7a77EdUwLbZElOI%2FuVnl5hZudqZ2gIVuydgaqbY807C4eX0tKBTXFVFeTh92NBsV%2B0%2FGRugBpc5lCbIz2q%2FNQP%2FrGMSWwtZKayYuNFTJXvYharMAbt9BnisKjHhrmJxYzCSq6LCBjqmAUbZiipVZQxfn9TGQh6ZdrbBwmc8TG0GdpzhcyFda0%2BVus5su3wiv%2BeZog%2Bk011GrSTBFp%2Bw0MoIFqKo3R3NQ3x9S%2FDRcOody464P9H2rEkdblOaLtV8uAbsFjG2KsHJ3xQP9%2Bzk5PWRYNqX79RRFIHiq7YLyHgIsT1huhyZL7lxwXFEugCI4XemsxRS1gESOkMSZ2Zl1vqp%2FMmaX8sfy0Z5HfvaB4Y%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2006%3A00%3A04%20GMT&X-Amz-Signature=5226ef5175fbe5bfa22dee6ce842fff8a46bef09524a98949a06b25c8860a2cd'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474758173/7Z3ova6Wtr-camera-video-segment-1749474760673.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSI2GHGTCK%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T220004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIE660DgdFCi8i8bvMivDROKxD%2ByCQRzirgbx79lzYgn3AiAyRf7fqefgxSZENBHtm9Aw%2FgI06qdyTqqpz4idAdpiLCrmAwjF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0qByG9j973ShKj46KroD61BMkswDVOuyWr%2F39rq%2F15a2nLb9X25cXbmqw525F0vEhltm3owMssOpJ8XG1g8fzPysBM%2FV8DLxsUIGLquTJnrJTQo5Tp3OpsWS2vuWz0LGciUFUlCibXMA10D0KSqlpjwHNY7hw5WwRRrZrkgY4NLes1xcPXZIl9BtLavqr0Xnb5VpSFVYIkfLaOeISCUymmd1F4PmsEPADaLjm9ZGSv9s5%2BupvKbNMkzTAevyWRTlcDwXQ7O8AkQNphRQy5kKvP2wGkM29t%2FzAFqXq85D%2FDgx%2F1bMS4PMbALbXh85jfqYC7pqCeu63rpwsHdsrywwsegg0QnXznhlAgFKzVtORnt6xSG6TPuvy3DP1wphmT7sXLxjx2kSk1wM9WdQBwr%2BbmZr10hR0dj0Tmv5kEXOVT4exoZK8VXihIepJBcSu7va%2BECIrOlyo1WtQE52mocX0TG45W%2B9MNB5J7a77EdUwLbZElOI%2FuVnl5hZudqZ2gIVuydgaqbY807C4eX0tKBTXFVFeTh92NBsV%2B0%2FGRugBpc5lCbIz2q%2FNQP%2FrGMSWwtZKayYuNFTJXvYharMAbt9BnisKjHhrmJxYzCSq6LCBjqmAUbZiipVZQxfn9TGQh6ZdrbBwmc8TG0GdpzhcyFda0%2BVus5su3wiv%2BeZog%2Bk011GrSTBFp%2Bw0MoIFqKo3R3NQ3x9S%2FDRcOody464P9H2rEkdblOaLtV8uAbsFjG2KsHJ3xQP9%2Bzk5PWRYNqX79RRFIHiq7YLyHgIsT1huhyZL7lxwXFEugCI4XemsxRS1gESOkMSZ2Zl1vqp%2FMmaX8sfy0Z5HfvaB4Y%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2006%3A00%3A04%20GMT&X-Amz-Signature=bad02bf708939c7e0fdf7e7df268c9978ca7f43790cf25e198ea2136035b1a5c'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474688613/JLWgs3viB6-camera-video-segment-driver-1749474691113.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSI2GHGTCK%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T220004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIE660DgdFCi8i8bvMivDROKxD%2ByCQRzirgbx79lzYgn3AiAyRf7fqefgxSZENBHtm9Aw%2FgI06qdyTqqpz4idAdpiLCrmAwjF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0qByG9j973ShKj46KroD61BMkswDVOuyWr%2F39rq%2F15a2nLb9X25cXbmqw525F0vEhltm3owMssOpJ8XG1g8fzPysBM%2FV8DLxsUIGLquTJnrJTQo5Tp3OpsWS2vuWz0LGciUFUlCibXMA10D0KSqlpjwHNY7hw5WwRRrZrkgY4NLes1xcPXZIl9BtLavqr0Xnb5VpSFVYIkfLaOeISCUymmd1F4PmsEPADaLjm9ZGSv9s5%2BupvKbNMkzTAevyWRTlcDwXQ7O8AkQNphRQy5kKvP2wGkM29t%2FzAFqXq85D%2FDgx%2F1bMS4PMbALbXh85jfqYC7pqCeu63rpwsHdsrywwsegg0QnXznhlAgFKzVtORnt6xSG6TPuvy3DP1wphmT7sXLxjx2kSk1wM9WdQBwr%2BbmZr10hR0dj0Tmv5kEXOVT4exoZK8VXihIepJBcSu7va%2BECIrOlyo1WtQE52mocX0TG45W%2B9MNB5J7a77EdUwLbZElOI%2FuVnl5hZudqZ2gIVuydgaqbY807C4eX0tKBTXFVFeTh92NBsV%2B0%2FGRugBpc5lCbIz2q%2FNQP%2FrGMSWwtZKayYuNFTJXvYharMAbt9BnisKjHhrmJxYzCSq6LCBjqmAUbZiipVZQxfn9TGQh6ZdrbBwmc8TG0GdpzhcyFda0%2BVus5su3wiv%2BeZog%2Bk011GrSTBFp%2Bw0MoIFqKo3R3NQ3x9S%2FDRcOody464P9H2rEkdblOaLtV8uAbsFjG2KsHJ3xQP9%2Bzk5PWRYNqX79RRFIHiq7YLyHgIsT1huhyZL7lxwXFEugCI4XemsxRS1gESOkMSZ2Zl1vqp%2FMmaX8sfy0Z5HfvaB4Y%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2006%3A00%3A04%20GMT&X-Amz-Signature=ad31a0d245695be2ea0e2a51e44fcd86764e43d94b518d6a3988b65254591c22'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474688613/LOHXs2uAgt-camera-video-segment-1749474691113.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSI2GHGTCK%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T220004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIE660DgdFCi8i8bvMivDROKxD%2ByCQRzirgbx79lzYgn3AiAyRf7fqefgxSZENBHtm9Aw%2FgI06qdyTqqpz4idAdpiLCrmAwjF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0qByG9j973ShKj46KroD61BMkswDVOuyWr%2F39rq%2F15a2nLb9X25cXbmqw525F0vEhltm3owMssOpJ8XG1g8fzPysBM%2FV8DLxsUIGLquTJnrJTQo5Tp3OpsWS2vuWz0LGciUFUlCibXMA10D0KSqlpjwHNY7hw5WwRRrZrkgY4NLes1xcPXZIl9BtLavqr0Xnb5VpSFVYIkfLaOeISCUymmd1F4PmsEPADaLjm9ZGSv9s5%2BupvKbNMkzTAevyWRTlcDwXQ7O8AkQNphRQy5kKvP2wGkM29t%2FzAFqXq85D%2FDgx%2F1bMS4PMbALbXh85jfqYC7pqCeu63rpwsHdsrywwsegg0QnXznhlAgFKzVtORnt6xSG6TPuvy3DP1wphmT7sXLxjx2kSk1wM9WdQBwr%2BbmZr10hR0dj0Tmv5kEXOVT4exoZK8VXihIepJBcSu7va%2BECIrOlyo1WtQE52mocX0TG45W%2B9MNB5J7a77EdUwLbZElOI%2FuVnl5hZudqZ2gIVuydgaqbY807C4eX0tKBTXFVFeTh92NBsV%2B0%2FGRugBpc5lCbIz2q%2FNQP%2FrGMSWwtZKayYuNFTJXvYharMAbt9BnisKjHhrmJxYzCSq6LCBjqmAUbZiipVZQxfn9TGQh6ZdrbBwmc8TG0GdpzhcyFda0%2BVus5su3wiv%2BeZog%2Bk011GrSTBFp%2Bw0MoIFqKo3R3NQ3x9S%2FDRcOody464P9H2rEkdblOaLtV8uAbsFjG2KsHJ3xQP9%2Bzk5PWRYNqX79RRFIHiq7YLyHgIsT1huhyZL7lxwXFEugCI4XemsxRS1gESOkMSZ2Zl1vqp%2FMmaX8sfy0Z5HfvaB4Y%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2006%3A00%3A04%20GMT&X-Amz-Signature=d05348b0e3bb17cf0dee126eed3d1cd1e32f1e5b1614bf50f97461565d2c7727'
